# ND01.xlsx realism pass — make the sample data look like genuine WMT
# extract output rather than placeholder 10/20/30… sequences, and tidy
# up the various Flag_* sheets so their single summary row reflects the
# correct flag code. Also nudges each sheet's saved cursor position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) WMT_Extract (sheet 1) — tier-count grid for OM_Key 1001/1002/1003
#    The placeholder 10/20/30/... filler values become mostly-zero
#    counts with a few small non-zero tallies, matching a realistic
#    extract.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WMT_Extract")

    # Row 2 (OM_Key 1001) tier counts
    $ws.Range("L2").Value = 0
    $ws.Range("M2").Value = 0
    $ws.Range("N2").Value = 2
    $ws.Range("O2").Value = 0
    $ws.Range("P2").Value = 0
    $ws.Range("Q2").Value = 0
    $ws.Range("R2").Value = 3
    $ws.Range("S2").Value = 0
    $ws.Range("T2").Value = 0
    $ws.Range("U2").Value = 0
    $ws.Range("V2").Value = 0
    $ws.Range("W2").Value = 2
    $ws.Range("X2").Value = 0
    $ws.Range("Y2").Value = 0
    $ws.Range("Z2").Value = 3
    $ws.Range("AA2").Value = 0
    $ws.Range("AB2").Value = 0
    $ws.Range("AC2").Value = 0
    $ws.Range("AD2").Value = 0
    $ws.Range("AE2").Value = 0
    $ws.Range("AF2").Value = 0
    $ws.Range("AG2").Value = 0
    $ws.Range("AH2").Value = 1
    $ws.Range("AI2").Value = 0
    $ws.Range("AJ2").Value = 0
    $ws.Range("AK2").Value = 0
    $ws.Range("AL2").Value = 0
    $ws.Range("AM2").Value = 0
    $ws.Range("AN2").Value = 0
    # Row 3 (OM_Key 1002) tier counts
    $ws.Range("L3").Value = 0
    $ws.Range("M3").Value = 3
    $ws.Range("N3").Value = 0
    $ws.Range("O3").Value = 0
    $ws.Range("P3").Value = 0
    $ws.Range("Q3").Value = 0
    $ws.Range("R3").Value = 0
    $ws.Range("S3").Value = 0
    $ws.Range("T3").Value = 2
    $ws.Range("U3").Value = 0
    $ws.Range("V3").Value = 0
    $ws.Range("W3").Value = 0
    $ws.Range("X3").Value = 0
    $ws.Range("Y3").Value = 0
    $ws.Range("Z3").Value = 1
    $ws.Range("AA3").Value = 0
    $ws.Range("AB3").Value = 0
    $ws.Range("AC3").Value = 0
    $ws.Range("AD3").Value = 0
    $ws.Range("AE3").Value = 0
    $ws.Range("AF3").Value = 0
    $ws.Range("AG3").Value = 0
    $ws.Range("AH3").Value = 0
    $ws.Range("AI3").Value = 0
    $ws.Range("AJ3").Value = 0
    $ws.Range("AK3").Value = 4
    $ws.Range("AL3").Value = 0
    $ws.Range("AM3").Value = 0
    $ws.Range("AN3").Value = 0
    # Row 4 (OM_Key 1003) tier counts
    $ws.Range("L4").Value = 0
    $ws.Range("M4").Value = 1
    $ws.Range("N4").Value = 3
    $ws.Range("O4").Value = 0
    $ws.Range("P4").Value = 0
    $ws.Range("Q4").Value = 4
    $ws.Range("R4").Value = 0
    $ws.Range("S4").Value = 0
    $ws.Range("T4").Value = 0
    $ws.Range("U4").Value = 2
    $ws.Range("V4").Value = 0
    $ws.Range("W4").Value = 0
    $ws.Range("X4").Value = 4

# Row 4 (OM_Key 1003) also gets a real LDU/team name instead of the
# Kainos placeholder — this introduces three brand-new shared strings.
$ws.Range("D4").Value = "Jonahs LDU"
$ws.Range("E4").Value = "JLDU"
$ws.Range("G4").Value = "JWMT"

# Saved cursor: scrolled right to show the later tier columns, with
# AB2 as the active cell.
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 22
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AB2").Select()

# ---------------------------------------------------------------------
# 2) Flag_Warr_4_n (sheet 4) — summary row's Row_Type flips from the
#    generic "I" to "N" (matching rows 2-3 on this sheet).
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Flag_Warr_4_n")
$ws4.Range("A4").Value = "N"
$ws4.Activate()
$ws4.Range("E24").Select()

# ---------------------------------------------------------------------
# 3) Flag_Upw (sheet 5) — summary row's Row_Type flips to "U".
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Flag_Upw")
$ws5.Range("A4").Value = "U"
$ws5.Activate()
$ws5.Range("A5").Select()

# ---------------------------------------------------------------------
# 4) Flag_O_Due (sheet 6) — summary row's Row_Type flips to "O".
# ---------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Flag_O_Due")
$ws6.Range("A4").Value = "O"
$ws6.Activate()
$ws6.Range("A5").Select()

# Leave WMT_Extract as the sheet on top when the workbook is reopened.
$ws.Activate()
